# Remove the confidential analyte entry "IgGI1H4N4S1" (row A9) by deleting
# the whole row and shifting the rows below it up, then clear the now
# trailing "test" row (previously A16, now A15) so the sheet ends with an
# empty cell instead of the leftover test value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 9 ("IgGI1H4N4S1") - shifts rows 10:16 up to 9:15
$ws.Rows("9:9").Delete()

# Clear the contents of what is now the last row (previously "test")
$ws.Range("A15").ClearContents()

# Update the active selection to match the resulting workbook state
$null = $ws.Range("F16").Select()
